$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Unmerged & Merged Cell Issue": several cells lost their checkmark
# ("ü") value and formatting (likely left blank after an unmerge).
# Restore the checkmark text and copy the formatting from a sibling cell
# in the same row that already has the correct style, so the resulting
# cell style exactly matches the other checked cells.

# G16 should look like the other "ü" cells in row 16 (e.g. C16)
$ws.Range("C16").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value = [char]0x00FC

# G19 should look like the other "ü" cells in row 19 (e.g. C19)
$ws.Range("C19").Copy()
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("G19").Value = [char]0x00FC

# G20 should look like the other "ü" cells in row 20 (e.g. C20)
$ws.Range("C20").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("G20").Value = [char]0x00FC

# D27 should look like the other shaded "ü" cells in row 27 (e.g. F27)
$ws.Range("F27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = [char]0x00FC

$excel.CutCopyMode = 0
